# Auto-generated edit script applying numeric corrections to the
# "currentAveragePrice*" / "LevePrice*" / "LeveProfit*" columns (H-N)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets, per the scheduled
# profit-data refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 1998.5
$ws.Range("I12").Value = 1998.5
$ws.Range("K12").Value = 1998.5
$ws.Range("M12").Value = -1828.5
$ws.Range("H33").Value = 162.63637
$ws.Range("I33").Value = 148.9
$ws.Range("K33").Value = 148.9
$ws.Range("M33").Value = 80.09999999999999
$ws.Range("H75").Value = 9999
$ws.Range("J75").Value = 9999
$ws.Range("L75").Value = 9999
$ws.Range("N75").Value = -11871
$ws.Range("H78").Value = 9999
$ws.Range("J78").Value = 9999
$ws.Range("L78").Value = 29997
$ws.Range("N78").Value = -39357
$ws.Range("H93").Value = 33900
$ws.Range("J93").Value = 33900
$ws.Range("L93").Value = 33900
$ws.Range("N93").Value = -38892
$ws.Range("H99").Value = 899.5
$ws.Range("I99").Value = 899.5
$ws.Range("K99").Value = 2698.5
$ws.Range("M99").Value = -1200.5
$ws.Range("H129").Value = 4850
$ws.Range("I129").Value = 4133.3335
$ws.Range("J129").Value = 7000
$ws.Range("K129").Value = 12400.0005
$ws.Range("L129").Value = 21000
$ws.Range("M129").Value = -7400.000499999998
$ws.Range("N129").Value = -31000
$ws.Range("H137").Value = 5024.8335
$ws.Range("I137").Value = 9625
$ws.Range("J137").Value = 2724.75
$ws.Range("K137").Value = 28875
$ws.Range("L137").Value = 8174.25
$ws.Range("M137").Value = -26325
$ws.Range("N137").Value = -13274.25
$ws.Range("H138").Value = 3050.1765
$ws.Range("J138").Value = 4553.6665
$ws.Range("L138").Value = 13660.9995
$ws.Range("N138").Value = -23940.9995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 55556924
$ws.Range("I2").Value = 66668130
$ws.Range("J2").Value = 900
$ws.Range("K2").Value = 66668130
$ws.Range("L2").Value = 900
$ws.Range("M2").Value = -66668017
$ws.Range("N2").Value = -1126
$ws.Range("H32").Value = 205044.05
$ws.Range("I32").Value = 753.45654
$ws.Range("J32").Value = 3337499.8
$ws.Range("K32").Value = 753.45654
$ws.Range("L32").Value = 3337499.8
$ws.Range("M32").Value = -466.45654
$ws.Range("N32").Value = -3338073.8
$ws.Range("H61").Value = 3076.9211
$ws.Range("I61").Value = 3073.4849
$ws.Range("K61").Value = 3073.4849
$ws.Range("M61").Value = -2861.4849
$ws.Range("H116").Value = 55556924
$ws.Range("I116").Value = 66668130
$ws.Range("J116").Value = 900
$ws.Range("K116").Value = 66668130
$ws.Range("L116").Value = 900
$ws.Range("M116").Value = -66665836
$ws.Range("N116").Value = -5488
$ws.Range("H136").Value = 3076.9211
$ws.Range("I136").Value = 3073.4849
$ws.Range("K136").Value = 9220.4547
$ws.Range("M136").Value = -6670.4547

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 55556924
$ws.Range("I3").Value = 66668130
$ws.Range("J3").Value = 900
$ws.Range("K3").Value = 66668130
$ws.Range("L3").Value = 900
$ws.Range("M3").Value = -66668016
$ws.Range("N3").Value = -1128
$ws.Range("H20").Value = 4398.2666
$ws.Range("I20").Value = 2797.9167
$ws.Range("K20").Value = 2797.9167
$ws.Range("M20").Value = -2550.9167
$ws.Range("H105").Value = 5213014
$ws.Range("I105").Value = 6948102.5
$ws.Range("K105").Value = 6948102.5
$ws.Range("M105").Value = -6946355.5
$ws.Range("H107").Value = 1317.1666
$ws.Range("I107").Value = 1317.1666
$ws.Range("K107").Value = 1317.1666
$ws.Range("M107").Value = 602.8334
$ws.Range("H134").Value = 4670.875
$ws.Range("I134").Value = 4385.3335
$ws.Range("K134").Value = 13156.0005
$ws.Range("M134").Value = -10621.0005

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2831.4243
$ws.Range("I58").Value = 2129.3845
$ws.Range("K58").Value = 2129.3845
$ws.Range("M58").Value = -1926.3845
$ws.Range("H132").Value = 2441.1
$ws.Range("I132").Value = 2323.4443
$ws.Range("J132").Value = 3500
$ws.Range("K132").Value = 6970.3329
$ws.Range("L132").Value = 10500
$ws.Range("M132").Value = -4440.3329
$ws.Range("N132").Value = -15560
$ws.Range("H134").Value = 3995.2727
$ws.Range("I134").Value = 4243.625
$ws.Range("K134").Value = 12730.875
$ws.Range("M134").Value = -10195.875
$ws.Range("H136").Value = 2831.4243
$ws.Range("I136").Value = 2129.3845
$ws.Range("K136").Value = 6388.1535
$ws.Range("M136").Value = -3838.1535

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 465.13043
$ws.Range("J5").Value = 454.66666
$ws.Range("L5").Value = 1363.99998
$ws.Range("N5").Value = -1587.99998
$ws.Range("H68").Value = 1530.9667
$ws.Range("I68").Value = 1471.5
$ws.Range("J68").Value = 1552.591
$ws.Range("K68").Value = 4414.5
$ws.Range("L68").Value = 4657.772999999999
$ws.Range("M68").Value = -3603.5
$ws.Range("N68").Value = -6279.772999999999
$ws.Range("H71").Value = 1530.9667
$ws.Range("I71").Value = 1471.5
$ws.Range("J71").Value = 1552.591
$ws.Range("K71").Value = 13243.5
$ws.Range("L71").Value = 13973.319
$ws.Range("M71").Value = -9187.5
$ws.Range("N71").Value = -22085.319
$ws.Range("H107").Value = 143726.28
$ws.Range("J107").Value = 167514
$ws.Range("L107").Value = 502542
$ws.Range("N107").Value = -506382
$ws.Range("H122").Value = 1573.5
$ws.Range("I122").Value = 1500
$ws.Range("J122").Value = 1598
$ws.Range("K122").Value = 13500
$ws.Range("L122").Value = 14382
$ws.Range("M122").Value = -11050
$ws.Range("N122").Value = -19282
$ws.Range("H132").Value = 3547.4
$ws.Range("I132").Value = 2637.8
$ws.Range("J132").Value = 4457
$ws.Range("K132").Value = 23740.2
$ws.Range("L132").Value = 40113
$ws.Range("M132").Value = -21210.2
$ws.Range("N132").Value = -45173
$ws.Range("H135").Value = 465.13043
$ws.Range("J135").Value = 454.66666
$ws.Range("L135").Value = 4091.99994
$ws.Range("N135").Value = -9161.99994

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H96").Value = 30261
$ws.Range("J96").Value = 30261
$ws.Range("L96").Value = 30261
$ws.Range("N96").Value = -35753
$ws.Range("H102").Value = 3641.7778
$ws.Range("I102").Value = 2642.625
$ws.Range("J102").Value = 4441.1
$ws.Range("K102").Value = 2642.625
$ws.Range("L102").Value = 4441.1
$ws.Range("M102").Value = -1020.625
$ws.Range("N102").Value = -7685.1
$ws.Range("H132").Value = 1747.421
$ws.Range("I132").Value = 1484.3077
$ws.Range("J132").Value = 2317.5
$ws.Range("K132").Value = 4452.9231
$ws.Range("L132").Value = 6952.5
$ws.Range("M132").Value = -1922.9231
$ws.Range("N132").Value = -12012.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H14").Value = 1000
$ws.Range("I14").Value = 1000
$ws.Range("K14").Value = 1000
$ws.Range("M14").Value = -828
$ws.Range("H22").Value = 1029.6364
$ws.Range("J22").Value = 1254.5
$ws.Range("L22").Value = 1254.5
$ws.Range("N22").Value = -1844.5
$ws.Range("H27").Value = 1029.6364
$ws.Range("J27").Value = 1254.5
$ws.Range("L27").Value = 1254.5
$ws.Range("N27").Value = -1468.5
$ws.Range("H46").Value = 4080
$ws.Range("I46").Value = 800
$ws.Range("K46").Value = 800
$ws.Range("M46").Value = -612
$ws.Range("H82").Value = 53077.184
$ws.Range("I82").Value = 10210.923
$ws.Range("J82").Value = 114995.11
$ws.Range("K82").Value = 10210.923
$ws.Range("L82").Value = 114995.11
$ws.Range("M82").Value = -9849.923000000001
$ws.Range("N82").Value = -115717.11
$ws.Range("H85").Value = 53077.184
$ws.Range("I85").Value = 10210.923
$ws.Range("J85").Value = 114995.11
$ws.Range("K85").Value = 10210.923
$ws.Range("L85").Value = 114995.11
$ws.Range("M85").Value = -8962.923000000001
$ws.Range("N85").Value = -117491.11
$ws.Range("H122").Value = 9750
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").Value = ""
$ws.Range("H132").Value = 4576.5386
$ws.Range("I132").Value = 5000
$ws.Range("K132").Value = 15000
$ws.Range("M132").Value = -12470

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H76").Value = 45000
$ws.Range("J76").Value = 45000
$ws.Range("L76").Value = 45000
$ws.Range("N76").Value = -45630
$ws.Range("H79").Value = 45000
$ws.Range("J79").Value = 45000
$ws.Range("L79").Value = 45000
$ws.Range("N79").Value = -47184
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").Value = ""
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").Value = ""
$ws.Range("H122").Value = 4273.3
$ws.Range("I122").Value = 4963.1665
$ws.Range("K122").Value = 14889.4995
$ws.Range("M122").Value = -12439.4995
$ws.Range("H124").Value = 100000
$ws.Range("J124").Value = 100000
$ws.Range("L124").Value = 100000
$ws.Range("N124").Value = -109820
$ws.Range("H132").Value = 6856.278
$ws.Range("I132").Value = 4207.7856
$ws.Range("J132").Value = 16126
$ws.Range("K132").Value = 12623.3568
$ws.Range("L132").Value = 48378
$ws.Range("M132").Value = -10093.3568
$ws.Range("N132").Value = -53438
$ws.Range("H136").Value = 4807.2085
$ws.Range("I136").Value = 4494.478
$ws.Range("J136").Value = 12000
$ws.Range("K136").Value = 13483.434
$ws.Range("L136").Value = 36000
$ws.Range("M136").Value = -10933.434
$ws.Range("N136").Value = -41100
